$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Cells.Item(2, 7).Value = 'servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Cells.Item(6, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Cells.Item(20, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(21, 7).Value = 'servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'norhan.mohamed@med.asu.edu.eg, yasmintarek@med.asu.edu.eg'
$ws.Cells.Item(28, 7).Value = 'yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(38, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(39, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(43, 7).Value = 'Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(44, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(47, 7).Value = 'merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(57, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(58, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(59, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(66, 7).Value = 'yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Cells.Item(77, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(78, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(81, 7).Value = 'NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(82, 7).Value = 'norhan.mohamed@med.asu.edu.eg, yasmintarek@med.asu.edu.eg'
$ws.Cells.Item(85, 7).Value = 'yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg'
$ws.Cells.Item(95, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(97, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(100, 7).Value = 'Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(104, 7).Value = 'merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(115, 7).Value = 'marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
